$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column H, matching the formatting already used by the
# other header cells in row 1 (bold / centered / bordered style).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Values for new "Save" column (1 = saved, 0 = blown/not saved)
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
